# Change 3rd party > 3rd
# Rename the "3rd Party Name" column header (cell B1 on Sheet1) to "3rd Name".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("B1").Value = "3rd Name"
